# Auto-generated edit script: applies numeric updates to match target diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1019.4
$ws.Range("I40").Value = 999
$ws.Range("K40").Value = 999
$ws.Range("M40").Value = -824
$ws.Range("H62").Value = 16543.436
$ws.Range("I62").Value = 17500
$ws.Range("J62").Value = 11999.75
$ws.Range("K62").Value = 17500
$ws.Range("L62").Value = 11999.75
$ws.Range("M62").Value = -16876
$ws.Range("N62").Value = -13247.75
$ws.Range("H65").Value = 16543.436
$ws.Range("I65").Value = 17500
$ws.Range("J65").Value = 11999.75
$ws.Range("K65").Value = 87500
$ws.Range("L65").Value = 59998.75
$ws.Range("M65").Value = -84380
$ws.Range("N65").Value = -66238.75
$ws.Range("H80").Value = 14686565
$ws.Range("I80").Value = 11111740
$ws.Range("J80").Value = 22729922
$ws.Range("K80").Value = 33335220
$ws.Range("L80").Value = 68189766
$ws.Range("M80").Value = -33334222
$ws.Range("N80").Value = -68191762
$ws.Range("H83").Value = 14686565
$ws.Range("I83").Value = 11111740
$ws.Range("J83").Value = 22729922
$ws.Range("K83").Value = 100005660
$ws.Range("L83").Value = 204569298
$ws.Range("M83").Value = -100000668
$ws.Range("N83").Value = -204579282
$ws.Range("H86").Value = 10003089
$ws.Range("I86").Value = 25001474
$ws.Range("J86").Value = 4165.6665
$ws.Range("K86").Value = 25001474
$ws.Range("L86").Value = 4165.6665
$ws.Range("M86").Value = -25000351
$ws.Range("N86").Value = -6411.6665
$ws.Range("H89").Value = 10003089
$ws.Range("I89").Value = 25001474
$ws.Range("J89").Value = 4165.6665
$ws.Range("K89").Value = 125007370
$ws.Range("L89").Value = 20828.3325
$ws.Range("M89").Value = -125001754
$ws.Range("N89").Value = -32060.3325
$ws.Range("H138").Value = 5447442
$ws.Range("I138").Value = 1267.4333
$ws.Range("J138").Value = 8082687.5
$ws.Range("K138").Value = 3802.2999
$ws.Range("L138").Value = 24248062.5
$ws.Range("M138").Value = 1337.7001
$ws.Range("N138").Value = -24258342.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H17").Value = 9000
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H88").Value = 20108.182
$ws.Range("I88").Value = 2033.3334
$ws.Range("J88").Value = 26886.25
$ws.Range("K88").Value = 2033.3334
$ws.Range("L88").Value = 26886.25
$ws.Range("M88").Value = -1627.3334
$ws.Range("N88").Value = -27698.25
$ws.Range("H91").Value = 20108.182
$ws.Range("I91").Value = 2033.3334
$ws.Range("J91").Value = 26886.25
$ws.Range("K91").Value = 2033.3334
$ws.Range("L91").Value = 26886.25
$ws.Range("M91").Value = -629.3334
$ws.Range("N91").Value = -29694.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1991227.2
$ws.Range("I105").Value = 7961010
$ws.Range("J105").Value = 1299.6666
$ws.Range("K105").Value = 7961010
$ws.Range("L105").Value = 1299.6666
$ws.Range("M105").Value = -7959263
$ws.Range("N105").Value = -4793.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 35717110
$ws.Range("I62").Value = 38464310
$ws.Range("K62").Value = 38464310
$ws.Range("M62").Value = -38463686
$ws.Range("H65").Value = 35717110
$ws.Range("I65").Value = 38464310
$ws.Range("K65").Value = 192321550
$ws.Range("M65").Value = -192318430

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1592.2222
$ws.Range("I2").Value = 1318.2858
$ws.Range("J2").Value = 2551
$ws.Range("K2").Value = 7909.714800000001
$ws.Range("L2").Value = 15306
$ws.Range("M2").Value = -7796.714800000001
$ws.Range("N2").Value = -15532
$ws.Range("H5").Value = 5637.143
$ws.Range("I5").Value = 905.3333
$ws.Range("K5").Value = 2715.9999
$ws.Range("M5").Value = -2603.9999
$ws.Range("H122").Value = 410.15625
$ws.Range("J122").Value = 782.0909
$ws.Range("L122").Value = 7038.8181
$ws.Range("N122").Value = -11938.8181
$ws.Range("H135").Value = 5637.143
$ws.Range("I135").Value = 905.3333
$ws.Range("K135").Value = 8147.9997
$ws.Range("M135").Value = -5612.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2385473.5
$ws.Range("I70").Value = 3337619.8
$ws.Range("J70").Value = 5108
$ws.Range("K70").Value = 3337619.8
$ws.Range("L70").Value = 5108
$ws.Range("M70").Value = -3337349.8
$ws.Range("N70").Value = -5648
$ws.Range("H73").Value = 2385473.5
$ws.Range("I73").Value = 3337619.8
$ws.Range("J73").Value = 5108
$ws.Range("K73").Value = 3337619.8
$ws.Range("L73").Value = 5108
$ws.Range("M73").Value = -3336683.8
$ws.Range("N73").Value = -6980
$ws.Range("H80").Value = 2711.111
$ws.Range("I80").Value = 2888.889
$ws.Range("J80").Value = 2533.3333
$ws.Range("K80").Value = 2888.889
$ws.Range("L80").Value = 2533.3333
$ws.Range("M80").Value = -1890.889
$ws.Range("N80").Value = -4529.3333
$ws.Range("H83").Value = 2711.111
$ws.Range("I83").Value = 2888.889
$ws.Range("J83").Value = 2533.3333
$ws.Range("K83").Value = 14444.445
$ws.Range("L83").Value = 12666.6665
$ws.Range("M83").Value = -9452.445
$ws.Range("N83").Value = -22650.6665
$ws.Range("H122").Value = 2987.3447
$ws.Range("I122").Value = 2764.158
$ws.Range("J122").Value = 3411.4
$ws.Range("K122").Value = 8292.474
$ws.Range("L122").Value = 10234.2
$ws.Range("M122").Value = -5842.474
$ws.Range("N122").Value = -15134.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 440.6154
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 448
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 448
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1038
$ws.Range("H27").Value = 440.6154
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 448
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 448
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -662
$ws.Range("H132").Value = 310941.16
$ws.Range("I132").Value = 81982.32000000001
$ws.Range("J132").Value = 719796.2
$ws.Range("K132").Value = 245946.96
$ws.Range("L132").Value = 2159388.6
$ws.Range("M132").Value = -243416.96
$ws.Range("N132").Value = -2164448.6
$ws.Range("H136").Value = 437117.97
$ws.Range("I136").Value = 716068.8
$ws.Range("J136").Value = 3194.4443
$ws.Range("K136").Value = 2148206.4
$ws.Range("L136").Value = 9583.332900000001
$ws.Range("M136").Value = -2145656.4
$ws.Range("N136").Value = -14683.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9203.071
$ws.Range("I132").Value = 1307.5
$ws.Range("J132").Value = 19730.5
$ws.Range("K132").Value = 3922.5
$ws.Range("L132").Value = 59191.5
$ws.Range("M132").Value = -1392.5
$ws.Range("N132").Value = -64251.5
$ws.Range("H136").Value = 529372.2
$ws.Range("I136").Value = 2356.8235
$ws.Range("K136").Value = 7070.470499999999
$ws.Range("M136").Value = -4520.470499999999
